# Add the new "bat 95" mod row (row 89) to the "mods" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mods")

$row = 89

# New shared strings must be introduced in the same order as the source
# workbook (M, then B, then F, then I) so the rebuilt sharedStrings.xml
# matches byte-for-byte.
$ws.Cells.Item($row, 13).Value = "bats 95%"
$ws.Cells.Item($row, 1).Value = "<Definition>"
$ws.Cells.Item($row, 2).Value = "invasion_bat_95"
$ws.Cells.Item($row, 3).Value = "gameplay"
$ws.Cells.Item($row, 4).Value = "levelUp"
$ws.Cells.Item($row, 5).Value = "spawn_frequency"
$ws.Cells.Item($row, 6).Value = "BatBig_Flock;BatSmall_Flock"
$ws.Cells.Item($row, 7).Value = 95
$ws.Cells.Item($row, 8).Value = "TID_EDIBLE_BAT_PL"
$ws.Cells.Item($row, 9).Value = "TID_QUIP_DRG_KILL_ENT_BATSMALL_02"
$ws.Cells.Item($row, 10).Value = "TID_EDIBLE_BAT_PL"
$ws.Cells.Item($row, 11).Value = "icon_tournament_rule"

# Copy the formatting from the row above (same look as every other mod row).
$ws.Range("A88:L88").Copy()
$ws.Range("A89:L89").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Expand the table to include the new row.
$table = $ws.ListObjects.Item("Table13")
$table.Resize($ws.Range("A3:L89"))

# Match the saved view state from the authored workbook.
$ws.Application.ActiveWindow.ScrollRow = 79
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("I89").Select()
